$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 6; this shifts the existing
# rows 6, 7, 8 down to 7, 8, 9 (and carries their formatting along).
$ws.Rows(6).Insert()

# Fill in the new row 6 with the new weekly record.
$ws.Cells.Item(6, 1).Value = 11
$ws.Cells.Item(6, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(6, 3).Value = "Bíobío"
$ws.Cells.Item(6, 4).Value = 44483
$ws.Cells.Item(6, 5).Value = 8
$ws.Cells.Item(6, 6).Value = 100112022
$ws.Cells.Item(6, 7).Value = "Arveja Verde"
$ws.Cells.Item(6, 8).Value = "Perfection"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 220
$ws.Cells.Item(6, 11).Value = 19000
$ws.Cells.Item(6, 12).Value = 20000
$ws.Cells.Item(6, 13).Value = 19455
$ws.Cells.Item(6, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(6, 15).Value = "Región Metropolitana"
$ws.Cells.Item(6, 16).Value = 778
$ws.Cells.Item(6, 17).Value = 25
$ws.Cells.Item(6, 18).Value = "Hortaliza"
